# Daily attendance processing - 2025-10-15 05:45:52
# Re-orders "Recorded By" email lists, updates a newly-recorded session
# (row 5) and refreshes the derived summary statistics that depend on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder "Recorded By" email lists (no data added/removed, just order) ---
$ws.Range("G2").Value = "eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G3").Value = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G4").Value = "rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# --- Row 5: session moved from "Pending" to "Recorded" ---
# Switch the row's formatting from the "Pending" yellow to the "Recorded"
# green by copying the format already used by recorded rows (e.g. row 2),
# then fill in the newly-recorded attendance data.
$ws.Range("A2:I2").Copy()
$ws.Range("A5:I5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G5").Value = "nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H5").Value = "1/221"
$ws.Range("I5").Value = "Recorded"

# --- Refresh summary stats (Recorded/Pending session counts, coverage) ---
# Percentages are stored as literal text (not numeric percent cells) in the
# source workbook, so a leading apostrophe is used to force text entry and
# avoid Excel's automatic "26.8%" -> 0.268 percent-number conversion.
$ws.Range("L6").Value = 15
$ws.Range("L8").Value = 40
$ws.Range("L9").Value = "'26.8%"
$ws.Range("L10").Value = "'40.7%"

# --- More "Recorded By" reorders ---
$ws.Range("G12").Value = "System, salma.elgendy.std@med.asu.edu.eg"
$ws.Range("G13").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# --- HISTOLOGY Year3/C1 per-subject rollup (row 15) ---
$ws.Range("O15").Value = 8
$ws.Range("Q15").Value = 20
$ws.Range("R15").Value = "'28.6%"
$ws.Range("S15").Value = "'46.3%"

# --- HISTOLOGY Year3/C2 per-subject rollup (row 16) ---
$ws.Range("S16").Value = "'34.4%"

# --- More "Recorded By" reorders ---
$ws.Range("G24").Value = "Salma.hassan@med.asu.edu.eg, marina_atef@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G25").Value = "marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"

$ws.Range("G31").Value = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G32").Value = "rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# --- Row 33: new recorder added, attendance count updated ---
$ws.Range("G33").Value = "nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H33").Value = "19/246"

$ws.Range("G40").Value = "System, salma.elgendy.std@med.asu.edu.eg"
$ws.Range("G41").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

$ws.Range("G52").Value = "Salma.hassan@med.asu.edu.eg, marina_atef@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G53").Value = "marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
